# aanpassing voor 3 map files
# - add 3 new countries (Belarus, Ukraine, Moldova) with their capitals and data
# - replace the "n/a" placeholder text with "No data"
# - add a concatenation helper column (H) that joins country/capital/temperatures

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- 1. Three new rows of data ---------------------------------------------
$ws.Range("A44").Value = "Belarus"
$ws.Range("A45").Value = "Ukraine"
$ws.Range("A46").Value = "Moldova"

$ws.Range("B45").Value = "Kiev"
$ws.Range("B44").Value = "Minsk"
$ws.Range("B46").Value = "Chisinau"

$ws.Range("C44").Value = 35
$ws.Range("D44").Value = 32
$ws.Range("E44").Value = "No data"
$ws.Range("F44").Value = "No data"

$ws.Range("C45").Value = 38
$ws.Range("D45").Value = 36
$ws.Range("E45").Value = 29
$ws.Range("F45").Value = 30

$ws.Range("C46").Value = 37
$ws.Range("D46").Value = 38
$ws.Range("E46").Value = 37
$ws.Range("F46").Value = "No data"

# --- 2. "n/a" -> "No data" everywhere it already appears on the sheet ------
$ws.Cells.Replace("n/a", "No data", -4163) | Out-Null

# --- 3. Helper column H: concatenate A..F with ".." between each value -----
# H4 is entered on its own, H5:H46 are entered as one block so Excel stores
# them as a shared formula (matching how the author filled the column down).
$ws.Range("H4").Formula = "=A4&""..""&B4&""..""&C4&""..""&D4&""..""&E4&""..""&F4"
$ws.Range("H5:H46").Formula = "=A5&""..""&B5&""..""&C5&""..""&D5&""..""&E5&""..""&F5"

# --- 4. Column H width -------------------------------------------------------
$ws.Columns("H").ColumnWidth = 44.85

# --- 5. View state: scroll down and select the new helper column -----------
$ws.Range("A25").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H4:H46").Select()
